# Añadiendo ASINs al fichero de descarga
#
# The sheet gains a new first data column "asins" (inserted before the
# existing "Imagen" column, shifting B:I -> C:J) and a brand-new row 6
# with a 5th product (Kinderkraft "Yummy" highchair). ASIN values are
# populated for all rows, and row 3's product data is replaced with the
# Kinderkraft "FINI" highchair record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0) - these scraped price fields use it as
# padding around the numeric text, matching the existing "118,99<nbsp>"
# cell already present in the sheet.
$nbsp = [char]0x00A0

# --- 1. Insert the new "asins" column before the current column B -------
$ws.Columns("B:B").Insert()

# Copy the header style (bold, bordered, centered) from the neighbouring
# header cell onto the new header cell, then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "asins"

# The inserted column picked up the bordered/bold style for its data rows
# too (Excel inherits format from the column being pushed); clear that so
# the data cells stay unstyled like the rest of the data grid.
$ws.Range("B2:B5").ClearFormats()

# --- 2. Populate ASIN values for the existing rows -----------------------
$ws.Range("B2").Value = "B06WVXN3HQ"
$ws.Range("B3").Value = "B071LG3DDN"
$ws.Range("B4").Value = "B072LNHXLQ"
$ws.Range("B5").Value = "B0785GR43J"

# --- 3. Replace row 3's product data with the Kinderkraft FINI record ---
$ws.Range("C3").Value = "https://images-na.ssl-images-amazon.com/images/I/41QeREFPHIL._AC_.jpg"
$ws.Range("D3").Value = "Kk KinderKraft"
$ws.Range("E3").Value = "Kinderkraft Trona Bebé 2 en 1 FINI, Silla Infantil, Ajustable, Segura, Gris"
$ws.Range("F3").Value = "89,00" + $nbsp
$ws.Range("G3").Value = $nbsp + "99,00" + $nbsp
# "2.807" parses as a number under the invariant locale used here; force
# text so it round-trips as the literal string like the other text cells.
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "2.807"
$ws.Range("H3").ClearFormats()
$ws.Range("I3").Value = "4,6"
$ws.Range("J3").Value = "La bandeja ajustable en 3 niveles fácil Ita el envío al niño en la trona"

# --- 4. Add the brand-new row 6 (Kinderkraft Yummy highchair) -----------
# Copy column A's numeric-index style (bold, bordered) down onto the new
# row before writing its value.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "B07GDVLMT5"
$ws.Range("C6").Value = "https://images-na.ssl-images-amazon.com/images/I/31KInCKXJlL._AC_.jpg"
$ws.Range("D6").Value = "Kinderkraft"
$ws.Range("E6").Value = "Kinderkraft Trona Bebé Ajustable Yummy, Segura, Bandeja, Hasta 3 los Años, Rosa"
$ws.Range("F6").Value = "78,89" + $nbsp
$ws.Range("G6").Value = $nbsp + "84,90" + $nbsp
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "339"
$ws.Range("H6").ClearFormats()
$ws.Range("I6").Value = "4,7"
$ws.Range("J6").Value = "MULTIFUNCIONAL: YUMMY es una trona para bebé apta para los niños que pueden estar sentados solos (desde aprox. los 6 meses de vida) hasta que pesen un máximo de 15 kg (aprox. 3 años). Gracias a la opción de plegado a un tamaño compacto es muy fácil guardarla. También la podrás llevar contigo a las vacaciones"

Write-Host "Done"
